# 678 importing cells with formula (#1473)
# Adds a new "Sheet4" worksheet (used by the new unit test) at the end of the
# workbook, with a small sample data table, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet so it lands at the end
# of the tab strip (Sheet1, Sheet2, Sheet3, Sheet4).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet4"

# Header row
$ws4.Range("A1").Value = "header1"
$ws4.Range("B1").Value = "header2"
$ws4.Range("C1").Value = "header3"

# Data row 2
$ws4.Range("A2").Value = "Concat_test1"
$ws4.Range("B2").Value = "str2"
$ws4.Range("C2").Value = "str3"
$ws4.Range("E2").Value = "test1"

# Data row 3
$ws4.Range("A3").Value = "Concat_test2"
$ws4.Range("B3").Value = "str5"
$ws4.Range("C3").Value = "str6"
$ws4.Range("E3").Value = "test2"

# Match the author's saved selection/active-cell state on the new sheet.
$ws4.Range("K15").Select()
